$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "岩山科技"
$ws.Range("B2").Value = "利欧股份"
$ws.Range("C2").Value = "利欧股份"

$ws.Range("A3").Value = "中国卫通"
$ws.Range("B3").Value = "岩山科技"
$ws.Range("C3").Value = "金风科技"

$ws.Range("A4").Value = "雷科防务"
$ws.Range("B4").Value = "中国卫通"
$ws.Range("C4").Value = "中国卫通"

$ws.Range("A5").Value = "中国卫星"
$ws.Range("B5").Value = "金风科技"
$ws.Range("C5").Value = "雷科防务"

$ws.Range("A6").Value = "航天电子"
$ws.Range("B6").Value = "雷科防务"
$ws.Range("C6").Value = "航天发展"

$ws.Range("A7").Value = "利欧股份"
$ws.Range("B7").Value = "航天电子"
$ws.Range("C7").Value = "中国卫星"

$ws.Range("A8").Value = "航天发展"
$ws.Range("B8").Value = "航天发展"
$ws.Range("C8").Value = "岩山科技"

$ws.Range("A9").Value = "金风科技"
$ws.Range("B9").Value = "山子高科"
$ws.Range("C9").Value = "山子高科"

$ws.Range("A10").Value = "蓝色光标"
$ws.Range("B10").Value = "蓝色光标"
$ws.Range("C10").Value = "鲁信创投"

$ws.Range("A11").Value = "浙江世宝"
$ws.Range("B11").Value = "中国卫星"
$ws.Range("C11").Value = "航天电子"

$ws.Range("A12").Value = "创新医疗"
$ws.Range("B12").Value = "东方财富"
$ws.Range("C12").Value = "海格通信"

$ws.Range("A13").Value = "山子高科"
$ws.Range("B13").Value = "海格通信"
$ws.Range("C13").Value = "天际股份"

$ws.Range("A14").Value = "三博脑科"
$ws.Range("B14").Value = "创新医疗"
$ws.Range("C14").Value = "北斗星通"

$ws.Range("A15").Value = "海格通信"
$ws.Range("B15").Value = "北斗星通"
$ws.Range("C15").Value = "平潭发展"

$ws.Range("A16").Value = "北斗星通"
$ws.Range("B16").Value = "三博脑科"
$ws.Range("C16").Value = "国晟科技"

$ws.Range("A17").Value = "三花智控"
$ws.Range("B17").Value = "三花智控"
$ws.Range("C17").Value = "西部材料"

$ws.Range("A18").Value = "南京熊猫"
$ws.Range("B18").Value = "浙江世宝"
$ws.Range("C18").Value = "浙江世宝"

$ws.Range("A19").Value = "天际股份"
$ws.Range("B19").Value = "中超控股"
$ws.Range("C19").Value = "万向钱潮"

$ws.Range("A20").Value = "普利特"
$ws.Range("B20").Value = "天际股份"
$ws.Range("C20").Value = "神剑股份"

$ws.Range("A21").Value = "东方财富"
$ws.Range("B21").Value = "南京熊猫"
$ws.Range("C21").Value = "泰和科技"

